$d = $word.ActiveDocument

# Helper-ish pattern: for each edit we
#   1) compute old/new paragraph text (sans trailing paragraph mark)
#   2) delete the paragraph's text range and re-insert the new text as ONE run
#   3) toggle Bold on/off over each desired sub-range to force a run split
#      at that boundary while keeping the run's formatting (<w:rPr/>) empty,
#      matching the identical-formatting runs produced by real Word when a
#      user types into the middle of existing text.

$LDQ = [char]0x201C
$RDQ = [char]0x201D

# ---------------------------------------------------------------------
# Paragraph 2: "Create database named ...." -> "Create " / "a " / "database named ...."
# ---------------------------------------------------------------------
$p = $d.Paragraphs(2)
$full = $p.Range.Text
$oldText = $full.Substring(0, $full.Length - 1)
$newText = "Create a database named " + $LDQ + "elections" + $RDQ + " in phpMyAdmin."

$start = $p.Range.Start
$end = $p.Range.End - 1
$r = $d.Range($start, $end)
$r.Delete()
$ins = $d.Range($start, $start)
$ins.InsertAfter($newText)

$split1 = $start + 7   # after "Create "
$split2 = $start + 9   # after "a "
$newEnd = $start + $newText.Length

$rA = $d.Range($start, $split1)
$rA.Font.Bold = 1
$rA.Font.Bold = 0
$rB = $d.Range($split1, $split2)
$rB.Font.Bold = 1
$rB.Font.Bold = 0
$rC = $d.Range($split2, $newEnd)
$rC.Font.Bold = 1
$rC.Font.Bold = 0

# ---------------------------------------------------------------------
# Paragraph 3: "Create table called ...." -> "Create table called .... database " /
#              "using the SQL provided" / "."
# ---------------------------------------------------------------------
$p = $d.Paragraphs(3)
$full = $p.Range.Text
$oldText = $full.Substring(0, $full.Length - 1)
$run1 = "Create table called " + $LDQ + "voters" + $RDQ + " in that database "
$run2 = "using the SQL provided"
$run3 = "."
$newText = $run1 + $run2 + $run3

$start = $p.Range.Start
$end = $p.Range.End - 1
$r = $d.Range($start, $end)
$r.Delete()
$ins = $d.Range($start, $start)
$ins.InsertAfter($newText)

$split1 = $start + $run1.Length
$split2 = $start + $run1.Length + $run2.Length
$newEnd = $start + $newText.Length

$rA = $d.Range($start, $split1)
$rA.Font.Bold = 1
$rA.Font.Bold = 0
$rB = $d.Range($split1, $split2)
$rB.Font.Bold = 1
$rB.Font.Bold = 0
$rC = $d.Range($split2, $newEnd)
$rC.Font.Bold = 1
$rC.Font.Bold = 0

# ---------------------------------------------------------------------
# Paragraph 14: delete the "hrep (VARCHAR, ...)" bullet entirely.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(14)
$p.Range.Delete()

# ---------------------------------------------------------------------
# Paragraph 25 (was 26): merge the 4 runs of the "and save this Excel Sheet..."
# paragraph into a single run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(25)
$full = $p.Range.Text
$textOnly = $full.Substring(0, $full.Length - 1)
$start = $p.Range.Start
$end = $p.Range.End - 1
$r = $d.Range($start, $end)
$r.Delete()
$ins = $d.Range($start, $start)
$ins.InsertAfter($textOnly)
$newEnd = $start + $textOnly.Length
$rAll = $d.Range($start, $newEnd)
$rAll.Font.Bold = 1
$rAll.Font.Bold = 0

# ---------------------------------------------------------------------
# Paragraph 26 (was 27): merge the 2 runs of the "Go to .../populate.php..."
# paragraph into a single run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(26)
$full = $p.Range.Text
$textOnly = $full.Substring(0, $full.Length - 1)
$start = $p.Range.Start
$end = $p.Range.End - 1
$r = $d.Range($start, $end)
$r.Delete()
$ins = $d.Range($start, $start)
$ins.InsertAfter($textOnly)
$newEnd = $start + $textOnly.Length
$rAll = $d.Range($start, $newEnd)
$rAll.Font.Bold = 1
$rAll.Font.Bold = 0

# ---------------------------------------------------------------------
# Paragraph 27 (was 28): merge the 3 runs of the "For voting on each
# individual mobile portal..." paragraph into a single run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(27)
$full = $p.Range.Text
$textOnly = $full.Substring(0, $full.Length - 1)
$start = $p.Range.Start
$end = $p.Range.End - 1
$r = $d.Range($start, $end)
$r.Delete()
$ins = $d.Range($start, $start)
$ins.InsertAfter($textOnly)
$newEnd = $start + $textOnly.Length
$rAll = $d.Range($start, $newEnd)
$rAll.Font.Bold = 1
$rAll.Font.Bold = 0

Write-Output "done"
